$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. All values are written as literal
# text (matching the original inlineStr cell type) by temporarily forcing a
# text number format, then restoring the default "Normal" style so no stray
# formatting is left behind.
$updates = [ordered]@{
    'D2' = '26.980.28'
    'E2' = '  +0.20%  '
    'D3' = '1.561.13'
    'E3' = '  +0.55%  '
    'E4' = '  -0.13%  '
    'D5' = '207.37'
    'E5' = '  +0.32%  '
    'E6' = '  +0.40%  '
    'D8' = '22.15'
    'E8' = '  +1.98%  '
    'E9' = '  +0.14%  '
    'D10' = '0.0597'
    'E10' = '  +1.92%  '
    'E11' = '  +0.02%  '
    'D12' = '1.782.56'
    'E12' = '  +0.51%  '
    'D13' = '1.559.48'
    'E13' = '  +0.46%  '
    'E15' = '  +0.85%  '
    'D16' = '62.08'
    'E16' = '  +0.68%  '
    'D17' = '26.960.71'
    'E17' = '  +0.19%  '
    'D18' = '217.33'
    'E18' = '  +0.16%  '
    'E19' = '  +2.12%  '
    'E20' = '  +2.11%  '
    'E21' = '  -0.17%  '
    'D22' = '4.10'
    'E22' = '  +1.27%  '
    'D23' = '9.21'
    'E23' = '  -0.20%  '
    'E24' = '  -1.33%  '
    'D25' = '153.55'
    'E25' = '  -0.21%  '
    'E26' = '  +0.38%  '
    'D27' = '15.07'
    'E27' = '  +1.31%  '
    'E28' = '  +1.42%  '
    'E29' = '  -0.14%  '
    'D30' = '0.0470'
    'E30' = '  +0.71%  '
    'E31' = '  +1.73%  '
    'E32' = '  +0.61%  '
    'E33' = '  +3.65%  '
    'D34' = '1.422.12'
    'E34' = '  +0.09%  '
    'E35' = '  +2.87%  '
    'D36' = '1.04'
    'E36' = '  +9.00%  '
    'E37' = '  +1.16%  '
    'E38' = '  +0.71%  '
    'E39' = '  +2.11%  '
    'D40' = '0.809'
    'E41' = '  -0.22%  '
    'D42' = '5.71'
    'E42' = '  -0.02%  '
    'E43' = '  +3.06%  '
    'E44' = '  +2.01%  '
    'D45' = '64.84'
    'E45' = '  +1.79%  '
    'D46' = '1.75'
    'E46' = '  +0.51%  '
    'D47' = '1.695.56'
    'E47' = '  +0.45%  '
    'D48' = '87.35'
    'E48' = '  +1.32%  '
    'E49' = '  -0.20%  '
    'E50' = '  -0.54%  '
    'D51' = '0.0954'
    'E51' = '  -0.31%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
